$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.212.73'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.682.70'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +0.34%  '
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '216.20'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -0.48%  '
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5245'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.91%  '
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06361'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.72%  '
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '21.49'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -1.82%  '
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07632'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +1.64%  '
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.679.74'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.513'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.5756'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.000008329'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -1.93%  '
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '66.05'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +2.26%  '
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.263.39'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.14%  '
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '4.868'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.81%  '
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.37%  '
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '189.44'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -0.34%  '
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.232'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +0.83%  '
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '148.88'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +2.73%  '
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '7.782'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '15.73'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.06314'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -2.94%  '
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.376'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +0.95%  '
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.315'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.60%  '
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.570'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.681'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.74%  '
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.6113'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -1.22%  '
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.418'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +0.60%  '
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.757'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +1.05%  '
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '6.189'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -1.44%  '
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01618'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -0.20%  '
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.8881'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +1.77%  '
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.097.65'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -1.48%  '
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -0.40%  '
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '100.45'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.833.07'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +0.35%  '
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.00000000109'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -3.65%  '
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '57.25'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +0.61%  '
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.065'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -1.36%  '
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.05273'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.4285'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.020'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -0.88%  '
$c.Style = "Normal"
